$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (isSelected), shifting
# isSelected/bandwidth/transRate/uploadTime/totalTime one column to the right.
$ws.Columns.Item(9).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 9).Value = "mu"

# New "mu" values for each data row (column I).
$muValues = @{
    2  = 289995.7640593715
    3  = 299772.0294312312
    4  = 231319.6593865472
    5  = 329844.9637816724
    6  = 297568.2072746158
    7  = 322193.0421175224
    8  = 291467.7789956704
    9  = 281081.506313911
    10 = 303112.4913603787
    11 = 280182.9275175941
    12 = 278320.2746903379
    13 = 301839.2237540928
    14 = 260105.1119513728
    15 = 274969.7264519503
    16 = 283118.8845010055
    17 = 289531.7998579305
    18 = 320461.7649340595
    19 = 302673.9683892264
    20 = 289833.6560853741
    21 = 280427.2764290712
}

foreach ($row in $muValues.Keys) {
    $ws.Cells.Item($row, 9).Value = $muValues[$row]
}
